$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings
# (e.g. "1.00", "410.00", "0.0400") are kept as literal text instead of
# being reinterpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '63.805.73'
$ws.Range("E2").Value = '  -2.02%  '

# Row 3
$ws.Range("D3").Value = '3.346.31'
$ws.Range("E3").Value = '  -2.29%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").Value = '546.42'
$ws.Range("E5").Value = '  -0.27%  '

# Row 6
$ws.Range("D6").Value = '172.23'
$ws.Range("E6").Value = '  -3.28%  '

# Row 7
$ws.Range("E7").Value = '  -3.60%  '

# Row 8
$ws.Range("D8").Value = '3.335.93'
$ws.Range("E8").Value = '  -2.43%  '

# Row 9
$ws.Range("E9").Value = '  +0.01%  '

# Row 10
$ws.Range("D10").Value = '0.612'

# Row 11
$ws.Range("E11").Value = '  +1.53%  '

# Row 12
$ws.Range("D12").Value = '53.64'
$ws.Range("E12").Value = '  +0.62%  '

# Row 13
$ws.Range("E13").Value = '  -1.53%  '

# Row 14
$ws.Range("E14").Value = '  -2.56%  '

# Row 15
$ws.Range("D15").Value = '3.874.72'
$ws.Range("E15").Value = '  -2.44%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '17.94'
$ws.Range("E16").Value = '  -1.41%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.354.32'
$ws.Range("E17").Value = '  -2.01%  '

# Row 18
$ws.Range("E18").Value = '  -3.05%  '

# Row 19
$ws.Range("D19").Value = '11.71'
$ws.Range("E19").Value = '  -0.47%  '

# Row 20
$ws.Range("D20").Value = '63.728.44'
$ws.Range("E20").Value = '  -2.31%  '

# Row 21
$ws.Range("D21").Value = '0.976'
$ws.Range("E21").Value = '  -0.26%  '

# Row 22
$ws.Range("D22").Value = '410.00'
$ws.Range("E22").Value = '  -0.67%  '

# Row 23
$ws.Range("E23").Value = '  +1.00%  '

# Row 24
$ws.Range("D24").Value = '4.35'
$ws.Range("E24").Value = '  +5.97%  '

# Row 25
$ws.Range("D25").Value = '13.69'
$ws.Range("E25").Value = '  +12.80%  '

# Row 26
$ws.Range("D26").Value = '83.03'
$ws.Range("E26").Value = '  -2.05%  '

# Row 27
$ws.Range("E27").Value = '  -1.71%  '

# Row 28
$ws.Range("E28").Value = '  -3.73%  '

# Row 29
$ws.Range("D29").Value = '8.61'
$ws.Range("E29").Value = '  -2.54%  '

# Row 30
$ws.Range("D30").Value = '29.07'
$ws.Range("E30").Value = '  -1.89%  '

# Row 31
$ws.Range("E31").Value = '  -1.75%  '

# Row 32
$ws.Range("D32").Value = '11.38'
$ws.Range("E32").Value = '  -1.87%  '

# Row 33
$ws.Range("D33").Value = '568.82'
$ws.Range("E33").Value = '  -6.53%  '

# Row 34
$ws.Range("D34").Value = '0.106'
$ws.Range("E34").Value = '  -1.36%  '

# Row 35
$ws.Range("D35").Value = '57.95'
$ws.Range("E35").Value = '  -1.56%  '

# Row 36
$ws.Range("E36").Value = '  +1.21%  '

# Row 37
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.12%  '

# Row 38
$ws.Range("D38").Value = '35.12'
$ws.Range("E38").Value = '  -5.40%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '3.38'
$ws.Range("E39").Value = '  +2.37%  '

# Row 40
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0737'
$ws.Range("E40").Value = '  -4.40%  '

# Row 41
$ws.Range("E41").Value = '  -2.53%  '

# Row 42
$ws.Range("D42").Value = '3.134.64'
$ws.Range("E42").Value = '  -1.28%  '

# Row 43
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.21%  '

# Row 44
$ws.Range("E44").Value = '  +1.23%  '

# Row 45
$ws.Range("D45").Value = '3.23'
$ws.Range("E45").Value = '  -1.16%  '

# Row 46
$ws.Range("D46").Value = '0.0400'
$ws.Range("E46").Value = '  -1.71%  '

# Row 47
$ws.Range("E47").Value = '  -4.33%  '

# Row 48
$ws.Range("D48").Value = '2.59'
$ws.Range("E48").Value = '  -4.14%  '

# Row 49
$ws.Range("E49").Value = '  -2.45%  '

# Row 50
$ws.Range("D50").Value = '132.24'
$ws.Range("E50").Value = '  -4.07%  '

# Row 51
$ws.Range("D51").Value = '8.07'
$ws.Range("E51").Value = '  -3.17%  '
